$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $result = $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $result) {
        Write-Host "WARNING: replacement not found for: $old"
    }
}

# --- Phase 1: replace each original text with a unique placeholder ---
# (text moves around in a cycle, so placeholders avoid clobbering)
Replace-Text 'Desenvolver nos alunos as competências necessárias para aplicar conceitos e ferramentas de estatística em bioprocessos, com ênfase em planejamento e otimização de experimentos.' '__PLACEHOLDER_OBJ_PT__'
Replace-Text 'Develop in students the necessary skills to apply concepts and tools of statistics in bioprocesses, with an emphasis on experimental design and optimization.' '__PLACEHOLDER_OBJ_EN__'
Replace-Text '5817181 - Valdeir Arantes' '__PLACEHOLDER_DOCENTE__'
Replace-Text '1. Fundamentos de estatística aplicada; 2. Análise de sistemas de medição; 3. Análise de Variância; 4. Testes de comparações múltiplas; 5. Controle estatístico de processos; 6. Planejamento de Experimentos: planejamentos fatoriais, superfícies de resposta, planejamentos de mistura; 7. Aplicação de software estatístico e estratégia sequencial de planejamentos experimentais.' '__PLACEHOLDER_RESUMO_PT__'
Replace-Text '1. The role of statistics in Engineering; Fundamentals of applied statistics; Analysis of Variance; Multiple comparison tests; Experimental Design' '__PLACEHOLDER_RESUMO_EN__'
Replace-Text '1. O papel da estatística na Engenharia: métodos de coleta de dados 2. Fundamentos de estatística aplicada 3. Análise de Variância: análise de variância de um modelo 4. Testes de comparações múltiplas (Tukey, Hsu) 5. Planejamento de Experimentos: vantagens dos experimentos fatoriais em relação aos experimentos do tipo um fator por vez; varielaboração do planejamento fatorial Completo do tipo 2^k e fracionado, e superfície de resposta' '__PLACEHOLDER_PROGRAMA_PT__'
Replace-Text 'A avaliação será composta por provas, exercícios, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n.' '__PLACEHOLDER_METODO_VAL__'
Replace-Text 'MF≥ 5,0 para aprovação 5,0' '__PLACEHOLDER_CRITERIO_VAL__'
Replace-Text '(MF+RC)/2 ≥ 5,0 para aprovação, onde RC é uma prova de recuperação a ser aplicada.' '__PLACEHOLDER_NORMA_VAL__'
Replace-Text '1. BOX, G.E.P.; HUNTER, W.G.; HUNTER, J.S. Statistics for Experimenters: an introduction to designs, data analysis and model building. New York: John Wiley & Sons Inc., 1978.^l^l2. RODRIGUES, M. I. e IEMMA, A. F. Planejamento de experimentos e otimização de processos. Campinas: Cárita editora, 2009.^l^l3. Planejamento e otimização de Experimentos. Roy E. Bruns, Edit. UNICAMP, 1996' '__PLACEHOLDER_BIBLIO__'

# --- Phase 2: replace placeholders with final text ---
Replace-Text '__PLACEHOLDER_OBJ_PT__' '1. Fundamentos de estatística aplicada; 2. Análise de sistemas de medição; 3. Análise de Variância; 4. Testes de comparações múltiplas; 5. Controle estatístico de processos; 6. Planejamento de Experimentos: planejamentos fatoriais, superfícies de resposta, planejamentos de mistura; 7. Aplicação de software estatístico e estratégia sequencial de planejamentos experimentais.'
Replace-Text '__PLACEHOLDER_OBJ_EN__' '1. The role of statistics in Engineering; Fundamentals of applied statistics; Analysis of Variance; Multiple comparison tests; Experimental Design'
Replace-Text '__PLACEHOLDER_DOCENTE__' 'Desenvolver nos alunos as competências necessárias para aplicar conceitos e ferramentas de estatística em bioprocessos, com ênfase em planejamento e otimização de experimentos.'
Replace-Text '__PLACEHOLDER_RESUMO_PT__' '1. O papel da estatística na Engenharia: métodos de coleta de dados 2. Fundamentos de estatística aplicada 3. Análise de Variância: análise de variância de um modelo 4. Testes de comparações múltiplas (Tukey, Hsu) 5. Planejamento de Experimentos: vantagens dos experimentos fatoriais em relação aos experimentos do tipo um fator por vez; varielaboração do planejamento fatorial Completo do tipo 2^k e fracionado, e superfície de resposta'
Replace-Text '__PLACEHOLDER_RESUMO_EN__' 'Develop in students the necessary skills to apply concepts and tools of statistics in bioprocesses, with an emphasis on experimental design and optimization.'
Replace-Text '__PLACEHOLDER_PROGRAMA_PT__' 'A avaliação será composta por provas, exercícios, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n.'
Replace-Text '__PLACEHOLDER_METODO_VAL__' 'MF≥ 5,0 para aprovação 5,0'
Replace-Text '__PLACEHOLDER_CRITERIO_VAL__' '(MF+RC)/2 ≥ 5,0 para aprovação, onde RC é uma prova de recuperação a ser aplicada.'
Replace-Text '__PLACEHOLDER_NORMA_VAL__' '1. BOX, G.E.P.; HUNTER, W.G.; HUNTER, J.S. Statistics for Experimenters: an introduction to designs, data analysis and model building. New York: John Wiley & Sons Inc., 1978.^l^l2. RODRIGUES, M. I. e IEMMA, A. F. Planejamento de experimentos e otimização de processos. Campinas: Cárita editora, 2009.^l^l3. Planejamento e otimização de Experimentos. Roy E. Bruns, Edit. UNICAMP, 1996'
Replace-Text '__PLACEHOLDER_BIBLIO__' '5817181 - Valdeir Arantes'
